$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new question row (row 27) - a "Bellring" story problem about Ramu
$ws.Range("A27").Value = "Ramu has {x} apples"
$ws.Range("B27").Value = "Bellring"
$ws.Range("C27").Value = "x1:5"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = "{x}"
$ws.Range("F27").Value = 10

# Update the view so the newly added row area is visible / selected
$ws.Activate()
$ws.Range("A21").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C35").Select() | Out-Null
